$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7086029052734375
$ws.Range("B1").Value = 1.009019613265991
$ws.Range("C1").Value = 4.411684513092041
$ws.Range("D1").Value = 2.1846923828125
$ws.Range("E1").Value = 1.64983594417572
